$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Append the new log entry as row 42
$row = 42
$ws.Cells.Item($row, 1).Value = "Algemene vraag"
$ws.Cells.Item($row, 2).Value = "klachten@testbedrijf123.nl"
$ws.Cells.Item($row, 3).Value = "Kunnen jullie mij meer informatie sturen over jullie diensten?"
$ws.Cells.Item($row, 4).Value = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item($row, 5).Value = "Bedankt, we hebben dit doorgestuurd naar support@testbedrijf123.nl."
$ws.Cells.Item($row, 6).Value = "2025-08-14 22:15:18"
$ws.Cells.Item($row, 7).Value = "Nee"
$ws.Cells.Item($row, 8).Value = "Ja"
$ws.Cells.Item($row, 9).Value = "Nee"
$ws.Cells.Item($row, 10).Value = "Nee"

# Extend the conditional formatting ranges to cover the new row
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $ws.Range("$col" + "2:" + "$col" + "41")
    $newRange = $ws.Range("$col" + "2:" + "$col" + "42")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary count for the affected category
$dash.Cells.Item(2, 2).Value = 34
